$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "76.105.29"
$ws.Range("E2").Value = "  +1.76%  "

$ws.Range("D3").Value = "2.918.99"
$ws.Range("E3").Value = "  +3.95%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "203.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "596.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.62%  "

$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("E8").Value = "  +0.57%  "

$ws.Range("E9").Value = "  +3.15%  "

$ws.Range("D10").Value = "2.917.62"
$ws.Range("E10").Value = "  +4.04%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.436"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +17.77%  "

$ws.Range("E12").Value = "  +0.61%  "

$ws.Range("E13").Value = "  +0.65%  "

$ws.Range("D14").Value = "3.456.54"
$ws.Range("E14").Value = "  +3.80%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.20"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.63%  "

$ws.Range("D16").Value = "76.052.36"
$ws.Range("E16").Value = "  +1.82%  "

$ws.Range("E17").Value = "  +1.81%  "

$ws.Range("D18").Value = "2.915.93"
$ws.Range("E18").Value = "  +3.69%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.82%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "373.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.93%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.77%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.42%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.37"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.80%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.03%  "

$ws.Range("E27").Value = "  +2.71%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.61%  "

$ws.Range("E29").Value = "  +4.44%  "

$ws.Range("E30").Value = "  +0.23%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.35%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "500.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.54%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.77"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.57%  "

$ws.Range("E34").Value = "  +3.01%  "

$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("E36").Value = "  +1.73%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "163.89"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.35%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.110"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +28.87%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.61"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.26%  "

$ws.Range("E40").Value = "  -4.11%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.369"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.14%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "183.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.09%  "

$ws.Range("E43").Value = "  -0.02%  "

$ws.Range("E44").Value = "  +0.32%  "

$ws.Range("E45").Value = "  +0.38%  "

$ws.Range("E46").Value = "  +0.53%  "

$ws.Range("E47").Value = "  -1.56%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.99%  "

$ws.Range("E49").Value = "  -0.15%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.13%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.42%  "
